# Applies scheduled Bahamut_Profits price/profit refresh across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1319.3049
$ws.Range("I15").Value = 1319.3049
$ws.Range("K15").Value = 3957.9147
$ws.Range("M15").Value = -3788.9147

$ws.Range("H116").Value = 2454.7083
$ws.Range("I116").Value = 2827.4375
$ws.Range("J116").Value = 1709.25
$ws.Range("K116").Value = 2827.4375
$ws.Range("L116").Value = 1709.25
$ws.Range("M116").Value = 614.5625
$ws.Range("N116").Value = -8593.25

$ws.Range("H132").Value = 2738.9
$ws.Range("I132").Value = 2693.5789
$ws.Range("K132").Value = 8080.736699999999
$ws.Range("M132").Value = -5550.736699999999

$ws.Range("H137").Value = 8334794
$ws.Range("I137").Value = 1460.5143
$ws.Range("J137").Value = 20001460
$ws.Range("K137").Value = 4381.5429
$ws.Range("L137").Value = 60004380
$ws.Range("M137").Value = -1831.5429
$ws.Range("N137").Value = -60009480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2258.5
$ws.Range("I102").Value = 2363
$ws.Range("J102").Value = 900
$ws.Range("K102").Value = 2363
$ws.Range("L102").Value = 900
$ws.Range("M102").Value = -741
$ws.Range("N102").Value = -4144

$ws.Range("H132").Value = 1781.5778
$ws.Range("I132").Value = 1298.2903
$ws.Range("K132").Value = 3894.8709
$ws.Range("M132").Value = -1364.8709

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1702.762
$ws.Range("I86").Value = 1619.1111
$ws.Range("J86").Value = 1853.3334
$ws.Range("K86").Value = 1619.1111
$ws.Range("L86").Value = 1853.3334
$ws.Range("M86").Value = -496.1111000000001
$ws.Range("N86").Value = -4099.3334

$ws.Range("H89").Value = 1702.762
$ws.Range("I89").Value = 1619.1111
$ws.Range("J89").Value = 1853.3334
$ws.Range("K89").Value = 8095.5555
$ws.Range("L89").Value = 9266.666999999999
$ws.Range("M89").Value = -2479.5555
$ws.Range("N89").Value = -20498.667

$ws.Range("H105").Value = 2692.577
$ws.Range("I105").Value = 2667
$ws.Range("J105").Value = 2800
$ws.Range("K105").Value = 2667
$ws.Range("L105").Value = 2800
$ws.Range("M105").Value = -920
$ws.Range("N105").Value = -6294

$ws.Range("H107").Value = 8954.412
$ws.Range("I107").Value = 1183
$ws.Range("J107").Value = 23202
$ws.Range("K107").Value = 1183
$ws.Range("L107").Value = 23202
$ws.Range("M107").Value = 737
$ws.Range("N107").Value = -27042

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1870.9375
$ws.Range("I31").Value = 1729.8334
$ws.Range("J31").Value = 1955.6
$ws.Range("K31").Value = 1729.8334
$ws.Range("L31").Value = 1955.6
$ws.Range("M31").Value = -1434.8334
$ws.Range("N31").Value = -2545.6

$ws.Range("H34").Value = 1870.9375
$ws.Range("I34").Value = 1729.8334
$ws.Range("J34").Value = 1955.6
$ws.Range("K34").Value = 1729.8334
$ws.Range("L34").Value = 1955.6
$ws.Range("M34").Value = -1527.8334
$ws.Range("N34").Value = -2359.6

$ws.Range("H134").Value = 3127.652
$ws.Range("I134").Value = 2511.3333
$ws.Range("J134").Value = 3800
$ws.Range("K134").Value = 7533.999899999999
$ws.Range("L134").Value = 11400
$ws.Range("M134").Value = -4998.999899999999
$ws.Range("N134").Value = -16470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 69.5625
$ws.Range("I33").Value = 55.916668
$ws.Range("J33").Value = 110.5
$ws.Range("K33").Value = 335.500008
$ws.Range("L33").Value = 663
$ws.Range("M33").Value = -52.50000799999998
$ws.Range("N33").Value = -1229

$ws.Range("H68").Value = 1354.662
$ws.Range("J68").Value = 1896.1316
$ws.Range("L68").Value = 5688.3948
$ws.Range("N68").Value = -7310.3948

$ws.Range("H71").Value = 1354.662
$ws.Range("J71").Value = 1896.1316
$ws.Range("L71").Value = 17065.1844
$ws.Range("N71").Value = -25177.1844

$ws.Range("H129").Value = 1446.6842
$ws.Range("I129").Value = 987.6667
$ws.Range("J129").Value = 1859.8
$ws.Range("K129").Value = 2963.0001
$ws.Range("L129").Value = 5579.4
$ws.Range("M129").Value = 2036.9999
$ws.Range("N129").Value = -15579.4

$ws.Range("H132").Value = 928.5714
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 980
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 8820
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -13880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 34111.11
$ws.Range("J124").Value = 34111.11
$ws.Range("L124").Value = 34111.11
$ws.Range("N124").Value = -43931.11

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1327.5
$ws.Range("I7").Value = 1186.6666
$ws.Range("J7").Value = 1750
$ws.Range("K7").Value = 1186.6666
$ws.Range("L7").Value = 1750
$ws.Range("M7").Value = -1074.6666
$ws.Range("N7").Value = -1974

$ws.Range("H40").Value = 11851298
$ws.Range("I40").Value = 12728838
$ws.Range("K40").Value = 12728838
$ws.Range("M40").Value = -12728702

$ws.Range("H126").Value = 1327.5
$ws.Range("I126").Value = 1186.6666
$ws.Range("J126").Value = 1750
$ws.Range("K126").Value = 3559.9998
$ws.Range("L126").Value = 5250
$ws.Range("M126").Value = -1089.9998
$ws.Range("N126").Value = -10190

$ws.Range("H132").Value = 1896564
$ws.Range("I132").Value = 2606597.5
$ws.Range("J132").Value = 3141.4167
$ws.Range("K132").Value = 7819792.5
$ws.Range("L132").Value = 9424.250100000001
$ws.Range("M132").Value = -7817262.5
$ws.Range("N132").Value = -14484.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1056.8
$ws.Range("I122").Value = 780
$ws.Range("J122").Value = 1126
$ws.Range("K122").Value = 2340
$ws.Range("L122").Value = 3378
$ws.Range("M122").Value = 110
$ws.Range("N122").Value = -8278

$ws.Range("H125").Value = 53320
$ws.Range("J125").Value = 53320
$ws.Range("L125").Value = 53320
$ws.Range("N125").Value = -63160

$ws.Range("H132").Value = 1718.7587
$ws.Range("I132").Value = 906.05
$ws.Range("K132").Value = 2718.15
$ws.Range("M132").Value = -188.1499999999996

$ws.Range("H136").Value = 2646.2827
$ws.Range("I136").Value = 3192.1155
$ws.Range("J136").Value = 1936.7
$ws.Range("K136").Value = 9576.3465
$ws.Range("L136").Value = 5810.1
$ws.Range("M136").Value = -7026.3465
$ws.Range("N136").Value = -10910.1
